$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. For the "Direct Submission" rows whose Accession column (F) holds only a
#    bare comma-separated list of GenBank accession numbers, re-order that
#    list into ascending order (it was previously listed newest-first).
$accessionCells = @("F26", "F77", "F53", "F54", "F60", "F75", "F48", "F49", "F74", "F73", "F76", "F70", "F8")
foreach ($coord in $accessionCells) {
    $cell = $ws.Range($coord)
    $val = $cell.Value2
    $parts = $val -split ", "
    $sortedParts = $parts | Sort-Object
    $cell.Value2 = [string]::Join(", ", $sortedParts)
}

# 2. Re-sort the data rows (A2:I81) in ascending order by the RefID column (A),
#    keeping the header row (row 1) fixed.
$sortRange = $ws.Range("A2:I81")
$sortKey = $ws.Range("A2")
$sortRange.Sort($sortKey, 1)

# 3. Update the sheet selection to reflect the new active column (F2:F81).
$ws.Range("F2:F81").Select()
